$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4392.3555
$ws.Range("I40").Value = 3428.4348
$ws.Range("K40").Value = 3428.4348
$ws.Range("M40").Value = -3253.4348
$ws.Range("H69").Value = 7615.4546
$ws.Range("J69").Value = 7807.4375
$ws.Range("L69").Value = 23422.3125
$ws.Range("N69").Value = -25170.3125
$ws.Range("H72").Value = 7615.4546
$ws.Range("J72").Value = 7807.4375
$ws.Range("L72").Value = 70266.9375
$ws.Range("N72").Value = -79002.9375
$ws.Range("H99").Value = 357.16666
$ws.Range("I99").Value = 359.66666
$ws.Range("K99").Value = 1078.99998
$ws.Range("M99").Value = 419.0000199999999
$ws.Range("H113").Value = 6882
$ws.Range("I113").Value = 4916.8335
$ws.Range("J113").Value = 8566.429
$ws.Range("K113").Value = 4916.8335
$ws.Range("L113").Value = 8566.429
$ws.Range("M113").Value = -1662.8335
$ws.Range("N113").Value = -15074.429
$ws.Range("H118").Value = 661
$ws.Range("I118").Value = 638.5714
$ws.Range("J118").Value = 700.25
$ws.Range("K118").Value = 1915.7142
$ws.Range("L118").Value = 2100.75
$ws.Range("M118").Value = -258.7142000000001
$ws.Range("N118").Value = -5414.75
$ws.Range("H138").Value = 2797.3718
$ws.Range("I138").Value = 1462.48
$ws.Range("J138").Value = 3427.0378
$ws.Range("K138").Value = 4387.440000000001
$ws.Range("L138").Value = 10281.1134
$ws.Range("M138").Value = 752.5599999999995
$ws.Range("N138").Value = -20561.1134

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4091.4893
$ws.Range("I32").Value = 2913.9534
$ws.Range("K32").Value = 2913.9534
$ws.Range("M32").Value = -2626.9534
$ws.Range("H61").Value = 7050.8887
$ws.Range("I61").Value = 5806.375
$ws.Range("K61").Value = 5806.375
$ws.Range("M61").Value = -5594.375
$ws.Range("H74").Value = 11497021
$ws.Range("I74").Value = 12822022
$ws.Range("K74").Value = 12822022
$ws.Range("M74").Value = -12821148
$ws.Range("H77").Value = 11497021
$ws.Range("I77").Value = 12822022
$ws.Range("K77").Value = 64110110
$ws.Range("M77").Value = -64105742
$ws.Range("H97").Value = 1400.7084
$ws.Range("I97").Value = 674.5789
$ws.Range("K97").Value = 674.5789
$ws.Range("M97").Value = -178.5789
$ws.Range("H132").Value = 2465.6155
$ws.Range("I132").Value = 1703.0605
$ws.Range("K132").Value = 5109.181500000001
$ws.Range("M132").Value = -2579.181500000001
$ws.Range("H135").Value = 41619.625
$ws.Range("J135").Value = 41619.625
$ws.Range("L135").Value = 41619.625
$ws.Range("N135").Value = -51759.625
$ws.Range("H136").Value = 7050.8887
$ws.Range("I136").Value = 5806.375
$ws.Range("K136").Value = 17419.125
$ws.Range("M136").Value = -14869.125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5055.8887
$ws.Range("I20").Value = 4643.077
$ws.Range("J20").Value = 6129.2
$ws.Range("K20").Value = 4643.077
$ws.Range("L20").Value = 6129.2
$ws.Range("M20").Value = -4396.077
$ws.Range("N20").Value = -6623.2
$ws.Range("H99").Value = 3488.6667
$ws.Range("I99").Value = 3414
$ws.Range("K99").Value = 3414
$ws.Range("M99").Value = -1916
$ws.Range("H105").Value = 23212.076
$ws.Range("I105").Value = 26944.75
$ws.Range("K105").Value = 26944.75
$ws.Range("M105").Value = -25197.75
$ws.Range("H132").Value = 69425.71000000001
$ws.Range("J132").Value = 69425.71000000001
$ws.Range("L132").Value = 69425.71000000001
$ws.Range("N132").Value = -79545.71000000001
$ws.Range("H135").Value = 50331.832
$ws.Range("J135").Value = 50331.832
$ws.Range("L135").Value = 50331.832
$ws.Range("N135").Value = -60471.832
$ws.Range("H137").Value = 57234.125
$ws.Range("J137").Value = 57234.125
$ws.Range("L137").Value = 57234.125
$ws.Range("N137").Value = -67434.125
$ws.Range("H138").Value = 64994.1
$ws.Range("J138").Value = 64994.1
$ws.Range("L138").Value = 64994.1
$ws.Range("N138").Value = -75274.10000000001
$ws.Range("H140").Value = 52476.832
$ws.Range("J140").Value = 52476.832
$ws.Range("L140").Value = 52476.832
$ws.Range("N140").Value = -62836.832

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34554.09
$ws.Range("I31").Value = 1764.5714
$ws.Range("K31").Value = 1764.5714
$ws.Range("M31").Value = -1469.5714
$ws.Range("H34").Value = 34554.09
$ws.Range("I34").Value = 1764.5714
$ws.Range("K34").Value = 1764.5714
$ws.Range("M34").Value = -1562.5714
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H58").Value = 4016.6667
$ws.Range("I58").Value = 1748.7273
$ws.Range("K58").Value = 1748.7273
$ws.Range("M58").Value = -1545.7273
$ws.Range("H74").Value = 76996.336
$ws.Range("J74").Value = 76996.336
$ws.Range("L74").Value = 76996.336
$ws.Range("N74").Value = -78744.336
$ws.Range("H77").Value = 76996.336
$ws.Range("J77").Value = 76996.336
$ws.Range("L77").Value = 230989.008
$ws.Range("N77").Value = -239725.008
$ws.Range("H104").Value = 42999.668
$ws.Range("J104").Value = 42999.668
$ws.Range("L104").Value = 42999.668
$ws.Range("N104").Value = -48241.668
$ws.Range("H107").Value = 3016.182
$ws.Range("I107").Value = 1455
$ws.Range("J107").Value = 5748.25
$ws.Range("K107").Value = 1455
$ws.Range("L107").Value = 5748.25
$ws.Range("M107").Value = 465
$ws.Range("N107").Value = -9588.25
$ws.Range("H132").Value = 3453.3333
$ws.Range("I132").Value = 2538.7896
$ws.Range("K132").Value = 7616.3688
$ws.Range("M132").Value = -5086.3688
$ws.Range("H133").Value = 41418.7
$ws.Range("J133").Value = 41418.7
$ws.Range("L133").Value = 41418.7
$ws.Range("N133").Value = -46478.7
$ws.Range("H135").Value = 69292.47
$ws.Range("J135").Value = 69292.47
$ws.Range("L135").Value = 69292.47
$ws.Range("N135").Value = -79432.47
$ws.Range("H136").Value = 4016.6667
$ws.Range("I136").Value = 1748.7273
$ws.Range("K136").Value = 5246.1819
$ws.Range("M136").Value = -2696.1819
$ws.Range("H138").Value = 69385.89
$ws.Range("J138").Value = 69385.89
$ws.Range("L138").Value = 69385.89
$ws.Range("N138").Value = -79665.89
$ws.Range("H140").Value = 94959.8
$ws.Range("J140").Value = 94959.8
$ws.Range("L140").Value = 94959.8
$ws.Range("N140").Value = -105319.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 946.25
$ws.Range("I36").Value = 142.5
$ws.Range("K36").Value = 427.5
$ws.Range("M36").Value = -258.5
$ws.Range("H113").Value = 1419.591
$ws.Range("I113").Value = 1156.2222
$ws.Range("J113").Value = 1601.9231
$ws.Range("K113").Value = 3468.6666
$ws.Range("L113").Value = 4805.7693
$ws.Range("M113").Value = -1298.6666
$ws.Range("N113").Value = -9145.7693
$ws.Range("H116").Value = 3366.375
$ws.Range("I116").Value = 1749.5
$ws.Range("K116").Value = 5248.5
$ws.Range("M116").Value = -1806.5
$ws.Range("H121").Value = 960.5
$ws.Range("J121").Value = 972.6
$ws.Range("L121").Value = 2917.8
$ws.Range("N121").Value = -5537.8
$ws.Range("H131").Value = 17748620
$ws.Range("I131").Value = 31251052
$ws.Range("J131").Value = 13890781
$ws.Range("K131").Value = 93753156
$ws.Range("L131").Value = 41672343
$ws.Range("M131").Value = -93748116
$ws.Range("N131").Value = -41682423

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 699.25
$ws.Range("I2").Value = 201
$ws.Range("J2").Value = 814.2308
$ws.Range("K2").Value = 201
$ws.Range("L2").Value = 814.2308
$ws.Range("M2").Value = -88
$ws.Range("N2").Value = -1040.2308
$ws.Range("H70").Value = 13777.714
$ws.Range("I70").Value = 13979.2
$ws.Range("K70").Value = 13979.2
$ws.Range("M70").Value = -13709.2
$ws.Range("H73").Value = 13777.714
$ws.Range("I73").Value = 13979.2
$ws.Range("K73").Value = 13979.2
$ws.Range("M73").Value = -13043.2
$ws.Range("H97").Value = 1333.0857
$ws.Range("I97").Value = 1039.0358
$ws.Range("J97").Value = 2509.2856
$ws.Range("K97").Value = 1039.0358
$ws.Range("L97").Value = 2509.2856
$ws.Range("M97").Value = -543.0358000000001
$ws.Range("N97").Value = -3501.2856
$ws.Range("H132").Value = 2407.4644
$ws.Range("I132").Value = 1261.3125
$ws.Range("K132").Value = 3783.9375
$ws.Range("M132").Value = -1253.9375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 11029.36
$ws.Range("I100").Value = 6522.6665
$ws.Range("J100").Value = 13564.375
$ws.Range("K100").Value = 6522.6665
$ws.Range("L100").Value = 13564.375
$ws.Range("M100").Value = -5981.6665
$ws.Range("N100").Value = -14646.375

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 61145.223
$ws.Range("J135").Value = 61145.223
$ws.Range("L135").Value = 61145.223
$ws.Range("N135").Value = -71285.223
$ws.Range("H136").Value = 3844.7
$ws.Range("I136").Value = 2647.6
$ws.Range("K136").Value = 7942.799999999999
$ws.Range("M136").Value = -5392.799999999999
